# A new weekly price record was added to the "Camote" (sweet potato) data set
# for "Vega Modelo de Temuco". It belongs chronologically right before the
# existing row 61, so insert a new row there (shifting every following row
# down by one) and populate it with the new observation's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("61:61").Insert()

$ws.Range("A61").Value = 10
$ws.Range("B61").Value = "Vega Modelo de Temuco"
$ws.Range("C61").Value = "La Araucanía"
$ws.Range("D61").Value = 45012
$ws.Range("E61").Value = 9
$ws.Range("F61").Value = 100114002
$ws.Range("G61").Value = "Camote"
$ws.Range("H61").Value = "Sin especificar"
$ws.Range("I61").Value = "Primera"
$ws.Range("J61").Value = 40
$ws.Range("K61").Value = 26000
$ws.Range("L61").Value = 26000
$ws.Range("M61").Value = 26000
$ws.Range("N61").Value = "$/malla 20 kilos"
$ws.Range("O61").Value = "Perú"
$ws.Range("P61").Value = 1300
$ws.Range("Q61").Value = 20
$ws.Range("R61").Value = "Hortaliza"
